# Auto-generated edit script applying cached-value updates to Sheets/Ultima_Profits.xlsx
# Each worksheet here corresponds to one Table_<Name> (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
# Values below are plain cached numbers (no formulas in the source workbook),
# so we just overwrite cell values directly to match the target state.

$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 12454.363
$ws.Range("I106").Value = 15642.857
$ws.Range("J106").Value = 6874.5
$ws.Range("K106").Value = 15642.857
$ws.Range("L106").Value = 6874.5
$ws.Range("M106").Value = -15011.857
$ws.Range("N106").Value = -8136.5
$ws.Range("H113").Value = 2460.5
$ws.Range("I113").Value = 2138.125
$ws.Range("J113").Value = 3750
$ws.Range("K113").Value = 2138.125
$ws.Range("L113").Value = 3750
$ws.Range("M113").Value = 1115.875
$ws.Range("N113").Value = -10258
$ws.Range("H132").Value = 9622521
$ws.Range("I132").Value = 5360.0713
$ws.Range("J132").Value = 20842542
$ws.Range("K132").Value = 16080.2139
$ws.Range("L132").Value = 62527626
$ws.Range("M132").Value = -13550.2139
$ws.Range("N132").Value = -62532686
$ws.Range("H138").Value = 4168956.8
$ws.Range("I138").Value = 8548459
$ws.Range("J138").Value = 3088.6829
$ws.Range("K138").Value = 25645377
$ws.Range("L138").Value = 9266.048699999999
$ws.Range("M138").Value = -25640237
$ws.Range("N138").Value = -19546.0487
$ws.Range("H141").Value = 1465.6471
$ws.Range("I141").Value = 1465.6471
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 4396.9413
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 783.0587000000005
$ws.Range("N141").ClearContents()

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1895860.9
$ws.Range("I45").Value = 2526948.8
$ws.Range("J45").Value = 2597.3333
$ws.Range("K45").Value = 2526948.8
$ws.Range("L45").Value = 2597.3333
$ws.Range("M45").Value = -2526571.8
$ws.Range("N45").Value = -3351.3333
$ws.Range("H74").Value = 1207.9143
$ws.Range("I74").Value = 1154.5555
$ws.Range("J74").Value = 1388
$ws.Range("K74").Value = 1154.5555
$ws.Range("L74").Value = 1388
$ws.Range("M74").Value = -280.5554999999999
$ws.Range("N74").Value = -3136
$ws.Range("H77").Value = 1207.9143
$ws.Range("I77").Value = 1154.5555
$ws.Range("J77").Value = 1388
$ws.Range("K77").Value = 5772.7775
$ws.Range("L77").Value = 6940
$ws.Range("M77").Value = -1404.7775
$ws.Range("N77").Value = -15676
$ws.Range("H97").Value = 6760.6875
$ws.Range("I97").Value = 9373.637000000001
$ws.Range("K97").Value = 9373.637000000001
$ws.Range("M97").Value = -8877.637000000001
$ws.Range("H110").Value = 638
$ws.Range("I110").Value = 692.125
$ws.Range("J110").Value = 205
$ws.Range("K110").Value = 692.125
$ws.Range("L110").Value = 205
$ws.Range("M110").Value = 1352.875
$ws.Range("N110").Value = -4295
$ws.Range("H128").Value = 49800
$ws.Range("J128").Value = 49800
$ws.Range("L128").Value = 49800
$ws.Range("N128").Value = -59760
$ws.Range("H132").Value = 5954615
$ws.Range("I132").Value = 7577892
$ws.Range("K132").Value = 22733676
$ws.Range("M132").Value = -22731146

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 10815.167
$ws.Range("I75").Value = 2931
$ws.Range("J75").Value = 50236
$ws.Range("K75").Value = 2931
$ws.Range("L75").Value = 50236
$ws.Range("M75").Value = -1995
$ws.Range("N75").Value = -52108
$ws.Range("H78").Value = 10815.167
$ws.Range("I78").Value = 2931
$ws.Range("J78").Value = 50236
$ws.Range("K78").Value = 8793
$ws.Range("L78").Value = 150708
$ws.Range("M78").Value = -4113
$ws.Range("N78").Value = -160068
$ws.Range("H94").Value = 764.8148
$ws.Range("I94").Value = 541.5
$ws.Range("J94").Value = 1402.8572
$ws.Range("K94").Value = 541.5
$ws.Range("L94").Value = 1402.8572
$ws.Range("M94").Value = -90.5
$ws.Range("N94").Value = -2304.8572

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 8351.857
$ws.Range("I16").Value = 10861.3
$ws.Range("J16").Value = 2078.25
$ws.Range("K16").Value = 10861.3
$ws.Range("L16").Value = 2078.25
$ws.Range("M16").Value = -10574.3
$ws.Range("N16").Value = -2652.25
$ws.Range("H53").Value = 24793.6
$ws.Range("J53").Value = 24793.6
$ws.Range("L53").Value = 24793.6
$ws.Range("N53").Value = -26007.6
$ws.Range("H58").Value = 2357
$ws.Range("I58").Value = 900.9231
$ws.Range("J58").Value = 6142.8
$ws.Range("K58").Value = 900.9231
$ws.Range("L58").Value = 6142.8
$ws.Range("M58").Value = -697.9231
$ws.Range("N58").Value = -6548.8
$ws.Range("H111").Value = 40700
$ws.Range("J111").Value = 40700
$ws.Range("L111").Value = 40700
$ws.Range("N111").Value = -48880
$ws.Range("H113").Value = 8351.857
$ws.Range("I113").Value = 10861.3
$ws.Range("J113").Value = 2078.25
$ws.Range("K113").Value = 10861.3
$ws.Range("L113").Value = 2078.25
$ws.Range("M113").Value = -8691.299999999999
$ws.Range("N113").Value = -6418.25
$ws.Range("H136").Value = 2357
$ws.Range("I136").Value = 900.9231
$ws.Range("J136").Value = 6142.8
$ws.Range("K136").Value = 2702.7693
$ws.Range("L136").Value = 18428.4
$ws.Range("M136").Value = -152.7692999999999
$ws.Range("N136").Value = -23528.4

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5265.8887
$ws.Range("I3").Value = 3305.6428
$ws.Range("J3").Value = 7376.923
$ws.Range("K3").Value = 9916.928400000001
$ws.Range("L3").Value = 22130.769
$ws.Range("M3").Value = -9804.928400000001
$ws.Range("N3").Value = -22354.769
$ws.Range("H129").Value = 2490
$ws.Range("I129").Value = 776.5
$ws.Range("J129").Value = 3497.9412
$ws.Range("K129").Value = 2329.5
$ws.Range("L129").Value = 10493.8236
$ws.Range("M129").Value = 2670.5
$ws.Range("N129").Value = -20493.8236
$ws.Range("H131").Value = 1119.0317
$ws.Range("I131").Value = 746
$ws.Range("J131").Value = 1151.1897
$ws.Range("K131").Value = 2238
$ws.Range("L131").Value = 3453.5691
$ws.Range("M131").Value = 2802
$ws.Range("N131").Value = -13533.5691

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 3514.7334
$ws.Range("I99").Value = 3514.7334
$ws.Range("K99").Value = 3514.7334
$ws.Range("M99").Value = -1268.7334
$ws.Range("H113").Value = 143544.42
$ws.Range("I113").Value = 143544.42
$ws.Range("K113").Value = 143544.42
$ws.Range("M113").Value = -141374.42
$ws.Range("H117").Value = 59310
$ws.Range("J117").Value = 59310
$ws.Range("L117").Value = 59310
$ws.Range("N117").Value = -66194
$ws.Range("H132").Value = 8262.117
$ws.Range("I132").Value = 8716
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 26148
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -23618
$ws.Range("N132").Value = -8060

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 228.26923
$ws.Range("I55").Value = 187.77777
$ws.Range("J55").Value = 249.70589
$ws.Range("K55").Value = 187.77777
$ws.Range("L55").Value = 249.70589
$ws.Range("M55").Value = -14.77777
$ws.Range("N55").Value = -595.70589
$ws.Range("H61").Value = 2233.9285
$ws.Range("I61").Value = 2134
$ws.Range("K61").Value = 2134
$ws.Range("M61").Value = -1932
$ws.Range("H93").Value = 1400.3846
$ws.Range("I93").Value = 1467.2222
$ws.Range("J93").Value = 1250
$ws.Range("K93").Value = 1467.2222
$ws.Range("L93").Value = 1250
$ws.Range("M93").Value = -219.2221999999999
$ws.Range("N93").Value = -3746
$ws.Range("H113").Value = 2233.9285
$ws.Range("I113").Value = 2134
$ws.Range("K113").Value = 2134
$ws.Range("M113").Value = 36
$ws.Range("H122").Value = 7333.6523
$ws.Range("I122").Value = 6687.1
$ws.Range("J122").Value = 7831
$ws.Range("K122").Value = 20061.3
$ws.Range("L122").Value = 23493
$ws.Range("M122").Value = -17611.3
$ws.Range("N122").Value = -28393
$ws.Range("H132").Value = 12508264
$ws.Range("I132").Value = 5178.0356
$ws.Range("J132").Value = 41682132
$ws.Range("K132").Value = 15534.1068
$ws.Range("L132").Value = 125046396
$ws.Range("M132").Value = -13004.1068
$ws.Range("N132").Value = -125051456
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -55060

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 201.2
$ws.Range("I113").Value = 183.35294
$ws.Range("J113").Value = 239.125
$ws.Range("K113").Value = 550.05882
$ws.Range("L113").Value = 717.375
$ws.Range("M113").Value = 1619.94118
$ws.Range("N113").Value = -5057.375
$ws.Range("H122").Value = 2805.6316
$ws.Range("I122").Value = 2831.375
$ws.Range("K122").Value = 8494.125
$ws.Range("M122").Value = -6044.125

